{"js": "// 1) The \"_GoBack\" bookmark currently sits right after \"${data_plan} \" in the\n//    \"El plan de datos...\" paragraph. It needs to move to just after the new\n//    \"N\u00b0: \" text further down in the document. Delete it from its old spot\n//    first -- it gets re-inserted (at id 0 again) once we find the new spot.\ntry {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n} catch (e) {\n  // Bookmark might already be absent -- nothing to clean up in that case.\n}\n\n// 2) Update the cached result of the \"Fecha:\" TIME field from\n//    \"11 de mayo de 2021\" to \"16 de mayo de 2021\".\nconst dateResults = context.document.body.search(\"11 de mayo de 2021\", {\n  matchCase: true,\n  matchWholeWord: false\n});\ndateResults.load(\"items\");\nawait context.sync();\n\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"16 de mayo de 2021\", \"Replace\");\n  await context.sync();\n}\n\n// 3) Change \"N\u00b0: 503\" to \"N\u00b0: \" (drop the hard-coded \"503\") and place the\n//    \"_GoBack\" bookmark immediately after that text (collapsed, zero-length),\n//    right before the \"${number}\" merge placeholder that follows it.\nconst numResults = context.document.body.search(\"N\u00b0: 503\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nnumResults.load(\"items\");\nawait context.sync();\n\nif (numResults.items.length > 0) {\n  const numRange = numResults.items[0];\n  numRange.insertText(\"N\u00b0: \", \"Replace\");\n  await context.sync();\n\n  const endOfNum = numRange.getRange(\"End\");\n  endOfNum.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# 1) The \"_GoBack\" bookmark currently sits right after \"${data_plan} \" in the\n#    \"El plan de datos...\" paragraph. It needs to move to just after the new\n#    \"N\u00b0: \" text further down in the document. Delete it from its old spot\n#    first -- it gets re-added (as id 0 again) once we find the new spot.\n$d = $word.ActiveDocument\n\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Update the cached result of the \"Fecha:\" TIME field from\n#    \"11 de mayo de 2021\" to \"16 de mayo de 2021\".\n$dateRange = $d.Content\n$dateRange.Find.Text = \"11 de mayo de 2021\"\nif ($dateRange.Find.Execute()) {\n    $dateRange.Text = \"16 de mayo de 2021\"\n}\n\n# 3) Change \"N\u00b0: 503\" to \"N\u00b0: \" (drop the hard-coded \"503\") and place the\n#    \"_GoBack\" bookmark immediately after that text (collapsed, zero-length),\n#    right before the \"${number}\" merge placeholder that follows it.\n$numRange = $d.Content\n$numRange.Find.Text = \"N\u00b0: 503\"\nif ($numRange.Find.Execute()) {\n    $numRange.Text = \"N\u00b0: \"\n    $numRange.Collapse(0)\n    $d.Bookmarks.Add(\"_GoBack\", $numRange)\n}\n"}
